$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Trim the per-language subcode sheets (yucu1253 / tani1257) down
#    to just the first "_a_0001" / "_b_0001" rows.
# ------------------------------------------------------------------
$yucu = $wb.Worksheets.Item("yucu1253")
$yucu.Rows("8:11").Delete()
$yucu.Rows("3:6").Delete()

$tani = $wb.Worksheets.Item("tani1257")
$tani.Rows("8:11").Delete()
$tani.Rows("3:6").Delete()

# ------------------------------------------------------------------
# 2. Insert a new "remarks" sheet right after "references", modeled
#    on the "references" sheet's glottocode layout.
# ------------------------------------------------------------------
$references = $wb.Worksheets.Item("references")
$remarks = $wb.Worksheets.Add($null, $references)
$remarks.Name = "remarks"

$remarks.Range("A1").Value = "glottocode"
$remarks.Range("B1").Value = "var001_remark"
$remarks.Range("C1").Value = "var002_remark"
$remarks.Range("D1").Value = "var003_remark"
$remarks.Range("A1:D1").Font.Bold = $true
$remarks.Range("A1:D1").HorizontalAlignment = -4108

$remarks.Range("A2").Value = "yucu1253"
$remarks.Range("A3").Value = "tani1257"

# ------------------------------------------------------------------
# 3. Bump the package version recorded on the "readme" sheet.
# ------------------------------------------------------------------
$readme = $wb.Worksheets.Item("readme")
$readme.Range("B6").Value = "version 0.0.71"

# ------------------------------------------------------------------
# 4. Restore the original active tab (inserting a sheet moves focus).
# ------------------------------------------------------------------
$yucu.Activate()
